# Applies the "Penalty Reward System" week-shift + forecast-value edits
# described in the commit diff:
#   - Sheet "Forecast Comparison": each week row's Week_Start_Date (col B)
#     is advanced by one week, and every MyForecast value (col D) becomes 3.
#   - Sheet "Summary": several derived metrics are refreshed to match.
#
# NOTE: Plain `Range.Value = "<numeric-or-date-like string>"` gets
# auto-coerced by Excel into a real number/date (exactly like typing it
# interactively), which would also implicitly create a new NumberFormat
# style and diverge from the original "plain text" cell type used
# throughout this workbook. To store these as literal text (matching the
# source file's inlineStr cells) without perturbing any styles, we build
# the text as a formula-computed string in a scratch cell, then copy/
# paste-special just the *values* into the destination - this keeps the
# destination a plain text value cell with no formula and no style change.

$wb = $excel.ActiveWorkbook
$scratchSheet = $wb.Worksheets.Item(1)
$scratch = $scratchSheet.Range("ZZ1")

function Set-TextValue {
    param($range, [string]$text)

    $scratch.Formula = '="' + $text + '"'
    $scratch.Copy()
    $range.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)
    $scratch.Clear()
}

# --- Sheet 1: "Forecast Comparison" ---
$ws1 = $wb.Worksheets.Item("Forecast Comparison")

$weekStartDates = @{
    2  = "2025-01-12"
    3  = "2025-01-19"
    4  = "2025-01-26"
    5  = "2025-02-02"
    6  = "2025-02-09"
    7  = "2025-02-16"
    8  = "2025-02-23"
    9  = "2025-03-02"
    10 = "2025-03-09"
    11 = "2025-03-16"
    12 = "2025-03-23"
    13 = "2025-03-30"
    14 = "2025-04-06"
    15 = "2025-04-13"
    16 = "2025-04-20"
    17 = "2025-04-27"
}

foreach ($row in 2..17) {
    Set-TextValue $ws1.Range("B$row") $weekStartDates[$row]
    $ws1.Range("D$row").Value = 3
}

# --- Sheet 2: "Summary" ---
$ws2 = $wb.Worksheets.Item("Summary")

Set-TextValue $ws2.Range("B2") "2024-02-11 to 2025-01-05"
Set-TextValue $ws2.Range("B8") "117 units"
Set-TextValue $ws2.Range("B9") "47"
Set-TextValue $ws2.Range("B10") "24"
Set-TextValue $ws2.Range("B11") "12"
Set-TextValue $ws2.Range("B12") "3"
Set-TextValue $ws2.Range("B13") "2025-01-12"
Set-TextValue $ws2.Range("B14") "3"
Set-TextValue $ws2.Range("B15") "2025-04-13"
